# SwaadSutra_Daily_2026-01-27.xlsx — 2026-01-27T14:09:57.836Z
#
# Records the day's single order on "Daily Orders", refreshes the roll-up
# counters on "Summary", and records the per-item breakdown on
# "Items Breakdown".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Daily Orders
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Daily Orders")

$ordersHeaders = @("Order ID","Date","Customer","Flat No","Phone","Items","Total","Status","Payment","Collection Date","Collection Time","Notes","Cancel Reason","Feedback")
for ($i = 0; $i -lt $ordersHeaders.Length; $i++) {
    $wsOrders.Cells.Item(1, $i + 1).Value = $ordersHeaders[$i]
}

$wsOrders.Cells.Item(2, 1).Value = 31
$wsOrders.Cells.Item(2, 2).Value = "2026-01-27 14:09"
$wsOrders.Cells.Item(2, 3).Value = "Girija Lakade"
$wsOrders.Cells.Item(2, 4).Value = "A 1507"

# Phone number: force text storage so the leading/standalone digits are not
# coerced into a number (matches the source report's numberStoredAsText hint).
$wsOrders.Range("E2").NumberFormat = "@"
$wsOrders.Cells.Item(2, 5).Value = "74996684"

$wsOrders.Cells.Item(2, 6).Value = "Jawar Bhakari x4"
$wsOrders.Cells.Item(2, 7).Value = 80
$wsOrders.Cells.Item(2, 8).Value = "NEW"
$wsOrders.Cells.Item(2, 9).Value = "PENDING"

# Collection Date / Collection Time / Notes / Cancel Reason / Feedback are
# still blank for a brand-new order, but the report keeps the row fully
# populated with empty text cells (a leading apostrophe enters an explicit
# empty text value instead of clearing the cell entirely).
$wsOrders.Cells.Item(2, 10).Value = "'"
$wsOrders.Cells.Item(2, 11).Value = "'"
$wsOrders.Cells.Item(2, 12).Value = "'"
$wsOrders.Cells.Item(2, 13).Value = "'"
$wsOrders.Cells.Item(2, 14).Value = "'"

# ---------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(2, 1).Value = 1   # Total Orders
$wsSummary.Cells.Item(2, 2).Value = 1   # New
$wsSummary.Cells.Item(2, 7).Value = 80  # Total Revenue

# ---------------------------------------------------------------------
# Sheet 3: Items Breakdown
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items Breakdown")
$wsItems.Cells.Item(1, 1).Value = "Item"
$wsItems.Cells.Item(1, 2).Value = "Quantity Ordered"
$wsItems.Cells.Item(1, 3).Value = "Revenue"

$wsItems.Cells.Item(2, 1).Value = "Jawar Bhakari"
$wsItems.Cells.Item(2, 2).Value = 4
$wsItems.Cells.Item(2, 3).Value = 80
